$d = $word.ActiveDocument

$old = "During this assessment I will be doing a few things in order to get the data to a usable state and then ill be using that data "
$new = "During this assessment I will be doing a few things in order to get the data to a usable state and then I" + [char]0x2019 + "ll be using that data "

$d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
